# Checked all functionalities and solved issues:
#  - Remove the now-unused extra rows (Kgf2 / row4 and robo / row5)
#  - Overwrite row 2 data with placeholder "r" values (testing/cleanup pass)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 and 5 entirely (shrinks used range to A1:M3)
$ws.Rows("4:5").Delete()

# Overwrite row 2 cells with "r" placeholder text (H2 gets "rr")
$ws.Range("A2").Value = "r"
$ws.Range("B2").Value = "r"
$ws.Range("C2").Value = "r"
$ws.Range("D2").Value = "r"
$ws.Range("E2").Value = "r"
$ws.Range("F2").Value = "r"
$ws.Range("G2").Value = "r"
$ws.Range("H2").Value = "rr"
$ws.Range("I2").Value = "r"
$ws.Range("J2").Value = "r"
$ws.Range("K2").Value = "r"
$ws.Range("L2").Value = "r"
$ws.Range("M2").Value = "r"

Write-Host "Row cleanup complete"
